# EIA Table 5.3 monthly refresh: October 2016 -> November 2016
# Adds the new "November" monthly row, shifts the trailing rolling-12 /
# year-to-date block down by one row, refreshes those computed figures,
# and updates the title + "Rolling 12 Months..." label text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new row for "November" right above the old "Rolling 12
#     Months Ending in October" banner row (row 52), pushing everything
#     below it down by one.
$ws.Rows("52:52").Insert()

# The freshly-inserted row picks up a generic/unformatted style; copy the
# formatting down from the row above (the previous "October" data row) so
# it matches the rest of the monthly-data rows exactly. Do this BEFORE
# writing the real values, since pasting the source cells' contents would
# otherwise clobber them.
$ws.Range("A51").Copy()
$ws.Range("A52").PasteSpecial(-4122)
$ws.Range("B51:F51").Copy()
$ws.Range("B52:F52").PasteSpecial(-4122)

# Populate the new November data row.
$ws.Range("A52").Value = "November"
$ws.Range("B52").Value = 12.75
$ws.Range("C52").Value = 10.25
$ws.Range("D52").Value = 6.64
$ws.Range("E52").Value = 9.0399999999999991
$ws.Range("F52").Value = 10.1

# --- "Year to Date" section (now starting at row 53, label text
#     unchanged) — refresh its three annual figures through November.
$ws.Range("B54").Value = 12.55
$ws.Range("C54").Value = 10.78
$ws.Range("D54").Value = 7.14
$ws.Range("E54").Value = 10.44
$ws.Range("F54").Value = 10.47

$ws.Range("B55").Value = 12.68
$ws.Range("C55").Value = 10.68
$ws.Range("D55").Value = 6.95
$ws.Range("E55").Value = 10.11
$ws.Range("F55").Value = 10.45

$ws.Range("B56").Value = 12.58
$ws.Range("C56").Value = 10.39
$ws.Range("D56").Value = 6.76
$ws.Range("E56").Value = 9.49
$ws.Range("F56").Value = 10.3

# --- "Rolling 12 Months Ending in ..." section (now starting at row 57
#     after the insert) — update its label to November and refresh the
#     two annual figures for the new rolling window.
$ws.Range("A57").Value = "Rolling 12 Months Ending in November"

$ws.Range("B58").Value = 12.64
$ws.Range("C58").Value = 10.66
$ws.Range("D58").Value = 6.93
$ws.Range("E58").Value = 10.14
$ws.Range("F58").Value = 10.42

$ws.Range("B59").Value = 12.56
$ws.Range("C59").Value = 10.37
$ws.Range("D59").Value = 6.73
$ws.Range("E59").Value = 9.51
$ws.Range("F59").Value = 10.27

# --- Title banner: "... 2006 - October 2016 ..." -> "... November 2016 ..."
$ws.Range("A2").Value = "Total by End-Use Sector, 2006 - November 2016 (Cents per Kilowatthour)"
